$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.680.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.97%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.087.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '345.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5158'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.98%  '

$ws.Range("E8").Value = '  -2.86%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09171'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '51.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.86%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.175'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.090.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.67%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.204'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.66%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.739'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("E17").Value = '  -1.77%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.009'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06658'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.198'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.757.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.02%  '

$ws.Range("E25").Value = '  -3.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.337.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.23%  '

$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.90'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.70%  '

$ws.Range("E31").Value = '  -4.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1053'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.636'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.202'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.34%  '

$ws.Range("E35").Value = '  -1.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.132'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02572'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06725'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2277'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6855'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.291'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6661'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.301'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.623'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.217'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000336'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '81.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.22%  '

$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.168'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.15%  '
